$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 118 ("Usain Bolt" quote), shifting all subsequent rows up by one.
$ws.Rows("118").Delete()
